# Daily attendance processing - 2025-10-05 09:15:25
# Applies the recorded-attendance refresh to the session analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Class Statistics summary (K/L column block) ---
# NOTE: the percentage cells store plain text (e.g. "29.4%"), not a numeric
# percent value - a leading apostrophe keeps Excel from re-interpreting the
# string as a number while leaving the cell's visible formatting untouched.
$ws.Range("L6").Value = 45          # Recorded Sessions
$ws.Range("L8").Value = 99          # Pending Sessions
$ws.Range("L9").Value = "'29.4%"    # Coverage %

# --- Row 14 / Row 31 (PHARMACOLOGY session 2, A1 group) : recorder order ---
$ws.Range("G14").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G31").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"

# --- Row 17 (PHYSIOLOGY session 1, A1 group) : recorder order ---
$ws.Range("G17").Value = "nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# --- Row 34 (same session, mirrored block) : recorder order ---
$ws.Range("G34").Value = "nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# --- Group A3 summary row (K17:S17) : one more session recorded ---
$ws.Range("O17").Value = 7          # Recorded
$ws.Range("Q17").Value = 9          # Pending
$ws.Range("R17").Value = "'41.2%"   # Coverage %
$ws.Range("S17").Value = "'45.1%"   # Avg Attendance %

# --- Group B3 / B4 Avg Attendance % (recalculated) ---
$ws.Range("S21").Value = "'63.7%"
$ws.Range("S22").Value = "'35.0%"

# --- Row 35 / Row 68 : recorder order ---
$ws.Range("G35").Value = "neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G68").Value = "neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"

# --- Row 45 / Row 62 : recorder order ---
$ws.Range("G45").Value = "mohamed.saleem@med.asu.edu.eg, backup@backdoor.com, System, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G62").Value = "mohamed.saleem@med.asu.edu.eg, backup@backdoor.com, System, Rania.a.youssef@med.asu.edu.eg"

# --- Row 51 : recorder order ---
$ws.Range("G51").Value = "neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"

# --- Row 52 : PHYSIOLOGY session 2 (A3 group) just got recorded.
# Copy the "Recorded" (green) formatting from the row above (already-recorded
# row 51) onto row 52, then fill in the recorder list, student count and status.
$srcFmt = $ws.Range("A51:I51")
$dstFmt = $ws.Range("A52:I52")
$srcFmt.Copy()
$dstFmt.PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G52").Value = "abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
$ws.Range("H52").Value = "22/220"
$ws.Range("I52").Value = "Recorded"

# --- Row 69 : recorder list gained a name ---
$ws.Range("G69").Value = "abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"

# --- Row 72 : recorder order ---
$ws.Range("G72").Value = "wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# --- Row 76 : recorder order ---
$ws.Range("G76").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"

# --- Row 83 / Row 150 : recorder order ---
$ws.Range("G83").Value = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg"
$ws.Range("G150").Value = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg"

# --- Row 85 / Row 102 : recorder order ---
$ws.Range("G85").Value = "neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G102").Value = "neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"

# --- Row 98 : recorder order ---
$ws.Range("G98").Value = "user@user.com, nourhanmohamed@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"

# --- Row 99 / Row 149 : recorder order ---
$ws.Range("G99").Value = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G149").Value = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"

# --- Row 109 : newly recorded students pushed the attendance count up ---
$ws.Range("H109").Value = "133/224"

# --- Row 116 / Row 133 : recorder order ---
$ws.Range("G116").Value = "afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G133").Value = "afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"

# --- Row 119 / Row 136 : recorder order ---
$ws.Range("G119").Value = "neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G136").Value = "neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"

# --- Row 126 : newly recorded students pushed the attendance count up ---
$ws.Range("H126").Value = "92/226"
